# cv122011a.xlsx - correção nos dados e inicio da analise PNAD 2009
#
# The original sheet had two "section header" rows (row 5 = "situação do
# domicílio", row 8 = "grandes regiões e unidades da federação") that only
# carried a label in column A with no data in B:H, which pushed every
# numeric value one row below the label it actually belongs to. This
# removes those two stray header rows (shifting the data back up so each
# label lines up with its real values) and relabels the "total" column
# header that had lost its caption.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two blank "section header" rows so the B:H data re-aligns
# with the correct row labels. Delete the lower row first so the row
# number of the upper one doesn't shift before it is removed.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()

# The second header row's second column lost its caption ("unnamed:
# 1_level_1" placeholder) - restore it to "total".
$ws.Range("B2").Value = "total"
